$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company name in row 2
$ws.Range("B2").Value = "Công Ty Lienketso"

# Add new "Email" column header (same bold style as the other headers)
$ws.Range("G1").Value = "Email"
$ws.Range("G1").Font.Bold = $true

# Add the Email values as hyperlinks (mailto links), matching rows 2 & 3
$ws.Range("G2").Value = "thanhan1507@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:thanhan1507@gmail.com")

$ws.Range("G3").Value = "who@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:who@gmail.com")

# Match final selection / active cell
[void]$ws.Range("J3").Select()
